$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the "_GoBack" bookmark that currently sits right after the
#    "4101-1" run (it will be re-created further down, around the new
#    split point in the second date).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. First deadline date (after "ІПЗ-41 "): 25.09.2024 -> 25.10.2024
#    Plain in-place text fix-up, formatting of the run is untouched.
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Start = 0
$findRange.Find.Execute("25.09.2024", $true, $false, $false, $false, $false, $true, 1, $false, "25.10.2024", 1) | Out-Null

# ------------------------------------------------------------------
# 3. Second deadline date (after "ІПЗ-42 "): 25.09.2024 -> 25.10.2024,
#    but this time the run is split in two ("25.10" / ".2024") with a
#    freshly placed "_GoBack" bookmark sitting right at the split.
# ------------------------------------------------------------------
$text = $d.Content.Text
$secondIdx = $text.LastIndexOf("25.09.2024")
$dateRange = $d.Range($secondIdx, $secondIdx + "25.09.2024".Length)
$dateRange.Text = "25.10.2024"

$splitPos = $secondIdx + "25.10".Length
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
